# The Markdown source for this deck used an image reference with a title,
# e.g. `![alt text](lalune.jpg "title")`. Previously only the link/alt
# text ("lalune.jpg") was written into PowerPoint's picture description
# (descr attribute / Shape.AlternativeText). This now also includes the
# title, prefixed with "fig:" per the existing "fig:" convention, giving
# the picture's description as "fig:  lalune.jpg".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Picture 1")
$sh.AlternativeText = "fig:  lalune.jpg"
